$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 4505
$ws.Range("L3").Value = 4827
$ws.Range("G4").Value = 1507
$ws.Range("L4").Value = 1197
$ws.Range("L5").Value = 280
$ws.Range("L6").Value = 4123
$ws.Range("G7").Value = 24734
$ws.Range("L7").Value = 14932

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 38

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L6").Value = 58
$ws.Range("L7").Value = 173

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 288
$ws.Range("L3").Value = 332
$ws.Range("L4").Value = 74
$ws.Range("L7").Value = 997

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L2").Value = 104
$ws.Range("L3").Value = 128
$ws.Range("L6").Value = 76
$ws.Range("L7").Value = 328

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L3").Value = 232
$ws.Range("L7").Value = 681

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L3").Value = 63
$ws.Range("L7").Value = 201

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 165
$ws.Range("L3").Value = 191
$ws.Range("L4").Value = 31
$ws.Range("L7").Value = 556

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L3").Value = 88
$ws.Range("L7").Value = 284

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L2").Value = 69
$ws.Range("L3").Value = 106
$ws.Range("L7").Value = 255

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L5").Value = 55
$ws.Range("L6").Value = 115
$ws.Range("L7").Value = 494
$ws.Range("L8").Value = 997
$ws.Range("L9").Value = 90
$ws.Range("L11").Value = 239
$ws.Range("L12").Value = 34
$ws.Range("L19").Value = 411
$ws.Range("L20").Value = 381
$ws.Range("L22").Value = 44
$ws.Range("L23").Value = 160
$ws.Range("L29").Value = 829
$ws.Range("L33").Value = 681
$ws.Range("L34").Value = 87
$ws.Range("L36").Value = 188
$ws.Range("L37").Value = 556
$ws.Range("L41").Value = 66
$ws.Range("L42").Value = 484
$ws.Range("L45").Value = 26
$ws.Range("L47").Value = 106
$ws.Range("L48").Value = 195
$ws.Range("L49").Value = 77
$ws.Range("L50").Value = 75
$ws.Range("L52").Value = 305
$ws.Range("L53").Value = 173
$ws.Range("L54").Value = 313
$ws.Range("L55").Value = 143
$ws.Range("G63").Value = 309
$ws.Range("L65").Value = 284
$ws.Range("L67").Value = 512
$ws.Range("L69").Value = 38
$ws.Range("L71").Value = 39
$ws.Range("L72").Value = 59
$ws.Range("L73").Value = 116
$ws.Range("L75").Value = 55
$ws.Range("L76").Value = 233
$ws.Range("L77").Value = 102
$ws.Range("L79").Value = 393
$ws.Range("L83").Value = 328
$ws.Range("L84").Value = 144
$ws.Range("L85").Value = 766
$ws.Range("L88").Value = 163
$ws.Range("L89").Value = 210
$ws.Range("L92").Value = 41
$ws.Range("L94").Value = 189
$ws.Range("L95").Value = 201
$ws.Range("L96").Value = 164
$ws.Range("L97").Value = 127
$ws.Range("L98").Value = 84
$ws.Range("L99").Value = 255
$ws.Range("G101").Value = 24734
$ws.Range("L101").Value = 14932

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L3").Value = 200
$ws.Range("L4").Value = 37
$ws.Range("L7").Value = 512

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L2").Value = 46
$ws.Range("L7").Value = 144

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L3").Value = 75
$ws.Range("L4").Value = 27
$ws.Range("L6").Value = 153
$ws.Range("L7").Value = 313

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 246
$ws.Range("L3").Value = 313
$ws.Range("L5").Value = 14
$ws.Range("L6").Value = 217
$ws.Range("L7").Value = 829

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L3").Value = 47
$ws.Range("L7").Value = 195

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L6").Value = 120
$ws.Range("L7").Value = 411

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L6").Value = 109
$ws.Range("L7").Value = 233

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("L4").Value = 11
$ws.Range("L7").Value = 115

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 66

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 139
$ws.Range("L3").Value = 165
$ws.Range("L6").Value = 134
$ws.Range("L7").Value = 484

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L3").Value = 48
$ws.Range("L5").Value = 2
$ws.Range("L7").Value = 143

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L6").Value = 44
$ws.Range("L7").Value = 160

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L6").Value = 46
$ws.Range("L7").Value = 164

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L3").Value = 143
$ws.Range("L7").Value = 393

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L6").Value = 103
$ws.Range("L7").Value = 381

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L2").Value = 69
$ws.Range("L3").Value = 54
$ws.Range("L6").Value = 50
$ws.Range("L7").Value = 188

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L3").Value = 161
$ws.Range("L7").Value = 494

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("L2").Value = 28
$ws.Range("L7").Value = 87

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L3").Value = 42
$ws.Range("L6").Value = 78
$ws.Range("L7").Value = 189

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 106

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("L2").Value = 18
$ws.Range("L7").Value = 84

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("L2").Value = 27
$ws.Range("L7").Value = 75

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L4").Value = 18
$ws.Range("L7").Value = 239

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("L3").Value = 37
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("L4").Value = 10
$ws.Range("L7").Value = 116

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L3").Value = 29
$ws.Range("L7").Value = 127

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L2").Value = 45
$ws.Range("L7").Value = 163

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L2").Value = 58
$ws.Range("L7").Value = 210

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("L2").Value = 13
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 232
$ws.Range("L3").Value = 310
$ws.Range("L6").Value = 159
$ws.Range("L7").Value = 766

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 39

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 59

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("L5").Value = 6
$ws.Range("L7").Value = 102

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("L2").Value = 5
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L2").Value = 103
$ws.Range("L3").Value = 96
$ws.Range("L5").Value = 6
$ws.Range("L7").Value = 305

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 34
